# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Mon Sep 25 16:48:58 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text (mixes plain decimals with thousands-dotted
# values like "26.351.75"), so keep it text-formatted to avoid Excel
# auto-coercing plain-looking decimals (e.g. '210.14') into numbers.
$ws.Range("D2").Value = "26.351.75"
$ws.Range("D3").Value = "1.592.29"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("E4").Value = "  -0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.14"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  -1.97%  "
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.245"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.63"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "1.813.97"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.604.71"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.08"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "26.349.24"
$ws.Range("E17").Value = "  -1.30%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").Value = "  +4.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "212.34"
$ws.Range("E20").Value = "  +1.73%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.93"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.82"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0506"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("D34").Value = "1.297.60"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.613"
$ws.Range("E35").Value = "  +4.12%  "
$ws.Range("E36").Value = "  -1.68%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.11"
$ws.Range("E39").Value = "  -10.90%  "
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.83"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "1.727.01"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.48"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.52"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.43"
$ws.Range("E51").Value = "  -1.77%  "
